$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 89 (old rows 89-108
# shift down to become rows 91-110). Two single-row inserts at the same
# index correctly shift everything down by two rows.
$ws.Rows.Item(89).Insert()
$ws.Rows.Item(89).Insert()

# Populate the two newly-inserted rows (89 and 90) with the new weekly
# price records. The non-varying columns (A, B, C, E, F, G, H, N, O, Q, R)
# repeat the same constants used throughout this data block.

# Row 89
$ws.Cells.Item(89, 1).Value2 = 1
$ws.Cells.Item(89, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(89, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(89, 4).Value2 = 44855
$ws.Cells.Item(89, 5).Value2 = 15
$ws.Cells.Item(89, 6).Value2 = 100112021
$ws.Cells.Item(89, 7).Value2 = "Ají"
$ws.Cells.Item(89, 8).Value2 = "Inferno"
$ws.Cells.Item(89, 9).Value2 = "Primera"
$ws.Cells.Item(89, 10).Value2 = 300
$ws.Cells.Item(89, 11).Value2 = 15000
$ws.Cells.Item(89, 12).Value2 = 16000
$ws.Cells.Item(89, 13).Value2 = 15500
$ws.Cells.Item(89, 14).Value2 = "$/caja 15 kilos"
$ws.Cells.Item(89, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(89, 16).Value2 = 1033
$ws.Cells.Item(89, 17).Value2 = 15
$ws.Cells.Item(89, 18).Value2 = "Hortaliza"

# Row 90
$ws.Cells.Item(90, 1).Value2 = 1
$ws.Cells.Item(90, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(90, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(90, 4).Value2 = 44855
$ws.Cells.Item(90, 5).Value2 = 15
$ws.Cells.Item(90, 6).Value2 = 100112021
$ws.Cells.Item(90, 7).Value2 = "Ají"
$ws.Cells.Item(90, 8).Value2 = "Inferno"
$ws.Cells.Item(90, 9).Value2 = "Segunda"
$ws.Cells.Item(90, 10).Value2 = 160
$ws.Cells.Item(90, 11).Value2 = 13000
$ws.Cells.Item(90, 12).Value2 = 14000
$ws.Cells.Item(90, 13).Value2 = 13375
$ws.Cells.Item(90, 14).Value2 = "$/caja 15 kilos"
$ws.Cells.Item(90, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(90, 16).Value2 = 892
$ws.Cells.Item(90, 17).Value2 = 15
$ws.Cells.Item(90, 18).Value2 = "Hortaliza"

# Make sure the new D89/D90 cells keep the same date-formatted style as
# the rest of column D in this block (the row-insert already copies the
# style from the row above, but re-assigning Value2 resets it, so restore
# the number format explicitly).
$dateFormat = $ws.Cells.Item(91, 4).NumberFormat
$ws.Cells.Item(89, 4).NumberFormat = $dateFormat
$ws.Cells.Item(90, 4).NumberFormat = $dateFormat
